# Mise à jour du classement - 26.03.2025 à 13:00

$wb = $excel.ActiveWorkbook

# --- Sheet "leaderboard2": "Qui a attrapé le plus de Cobblemons ?" ---
$ws1 = $wb.Worksheets.Item("leaderboard2")
$ws1.Range("D3").Value = 57
$ws1.Range("D4").Value = 31
$ws1.Range("D5").Value = 17
$ws1.Range("B13").Value = "Dernière update le 26.03.25 à 13:00"

# --- Sheet "leaderboard3": "Qui a attrapé le plus de Shiny Cobblemons ?" ---
$ws2 = $wb.Worksheets.Item("leaderboard3")
$ws2.Range("C3").Value = "Lokys"
$ws2.Range("D3").Value = 6
$ws2.Range("C4").Value = "BKZRackham"
$ws2.Range("D4").Value = 6
$ws2.Range("C5").Value = "ArtyumsM"
$ws2.Range("D5").Value = 3
$ws2.Range("C6").Value = "Machoppeur_"
$ws2.Range("B13").Value = "Dernière update le 26.03.25 à 13:00"
